$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

$ws.Cells.Item($row, 1).Value = "RFTFFD"
$ws.Cells.Item($row, 2).Value = "Cabezal para dfx-8000/8500"
$ws.Cells.Item($row, 3).Value = "DFX 8000 8500"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 1200000
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E19-D19)*G19"
$ws.Cells.Item($row, 9).Formula = "=D19*F19"
$ws.Cells.Item($row, 10).Value = 0
